# MI4 - Logboek Eilish Van Der Snickt
# Commit: "FireStore compatibel maken met de verschillende users"
#
# Adds the "Week 12" total time, a new "Week 13" header/section with two
# new logboek entries (rows 37-38) and three extra "Interessante links"
# hyperlinks (rows 38-40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 12 total time (B33) -------------------------------------------
$ws.Range("B33").Value = "8 uur"
$ws.Range("B33").Font.Bold = $true

# --- Week 13 header (A36) -------------------------------------------------
$ws.Range("A36").Value = "Week 13"
$ws.Range("A36").Font.Bold = $true

# --- Row 37: new logboek entry ------------------------------------------
$ws.Range("A37").Value = 43592
$ws.Range("A37").NumberFormat = "d-mmm"
$ws.Range("B37").Value = "20 minuten"
$ws.Range("C37").Value = "locatie route testen"

# --- Row 38: new logboek entry ------------------------------------------
$ws.Range("A38").Value = 43595
$ws.Range("A38").NumberFormat = "d-mmm"
$ws.Range("B38").Value = "4 uur"
$ws.Range("C38").Value = "testen of de routes in de juiste documenten opgeslagen en opgehaald worden, RouteGegevens collection juist wegschrijven en ophalen, viewpager tabs opvullen met gegevens uit firestore. Firestore compatibel maken per gebruiker, wegschrijven en ophalen van de juiste info uit de firestore per gebruiker"

# --- Extra "Interessante links" hyperlinks (Q38:Q40) ---------------------
$ws.Range("Q38").Value = "http://www.downloadinformer.com/how-to-use-switch-statement-in-androidjava/"
$ws.Hyperlinks.Add($ws.Range("Q38"), "http://www.downloadinformer.com/how-to-use-switch-statement-in-androidjava/") | Out-Null
$ws.Range("Q38").Style = "Hyperlink"

$ws.Range("Q39").Value = "https://stackoverflow.com/questions/48873465/how-to-add-sub-collection-to-a-document-in-firestore"
$ws.Hyperlinks.Add($ws.Range("Q39"), "https://stackoverflow.com/questions/48873465/how-to-add-sub-collection-to-a-document-in-firestore") | Out-Null
$ws.Range("Q39").Style = "Hyperlink"

$ws.Range("Q40").Value = "https://angularfirebase.com/lessons/managing-firebase-user-relationships-to-database-records/"
$ws.Hyperlinks.Add($ws.Range("Q40"), "https://angularfirebase.com/lessons/managing-firebase-user-relationships-to-database-records/") | Out-Null
$ws.Range("Q40").Style = "Hyperlink"

# --- Update the visible selection / scroll position -----------------------
$ws.Range("C39").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
